$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprints")

# --- Fix existing Sprint 8 rows (content-only changes) ---
# Row 34: "QUEM REALIZOU" changes from "Cristiele e Gabriel" to "Cristielen"
$ws.Range("E34").Value = "Cristielen"

# --- Grow the Sprints table (Tabela1) by one row ---
$lo = $ws.ListObjects.Item("Tabela1")
$newListRow = $lo.ListRows.Add()

# Row 40 becomes the new blank separator row: copy formatting/style from the
# old blank separator row (39) down into the newly added row.
$ws.Range("B39:F39").Copy($ws.Range("B40:F40"))

# Row 39 becomes a populated data row: copy formatting/style from the row
# above it (38), then overwrite with the new row's values.
$ws.Range("B38:F38").Copy($ws.Range("B39:F39"))

$ws.Range("B39").Value = "Sprint 8"
$ws.Range("C39").Value = "08/05/2024 - 15/05/2024"
$ws.Range("D39").Value = "Alterações no Manual do Usuário"
$ws.Range("E39").Value = "Bruno"
$ws.Range("F39").Value = "Pronto"

# Row 35: "TAREFA" changes from "Edições nos Gráficos" to "Edições nos Diagramas"
$ws.Range("D35").Value = "Edições nos Diagramas"

# --- Update the view so the selection matches the edited area ---
$ws.Range("D35").Select()
